$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.72563878218237
$ws.Range("C2").Value = 5.355237703425621
$ws.Range("D2").Value = 5.988194326377732
$ws.Range("E2").Value = 16.41747679856375
$ws.Range("G2").Value = 3.641018890355125
$ws.Range("K2").Value = 9.970682741257116
$ws.Range("N2").Value = 18.65873008793067
$ws.Range("O2").Value = 23.82510249221289
$ws.Range("B3").Value = 10.41567965699306
$ws.Range("C3").Value = 5.115359297618998
$ws.Range("D3").Value = 5.870128634375859
$ws.Range("E3").Value = 15.49148942939535
$ws.Range("G3").Value = 3.643287013654936
$ws.Range("K3").Value = 9.753396103095719
$ws.Range("N3").Value = 18.72427094975153
$ws.Range("O3").Value = 23.8501096054148
$ws.Range("B4").Value = 10.22319461811687
$ws.Range("C4").Value = 4.961083630589762
$ws.Range("D4").Value = 5.79822248308855
$ws.Range("E4").Value = 14.89867955553206
$ws.Range("G4").Value = 3.64475197385754
$ws.Range("K4").Value = 9.619785751307392
$ws.Range("N4").Value = 18.76632860322066
$ws.Range("O4").Value = 23.87130860129606
$ws.Range("B5").Value = 10.14434075744446
$ws.Range("C5").Value = 4.896502571506255
$ws.Range("D5").Value = 5.769112731616776
$ws.Range("E5").Value = 14.65127914793525
$ws.Range("G5").Value = 3.645367206525984
$ws.Range("K5").Value = 9.565370976253037
$ws.Range("N5").Value = 18.78392530236685
$ws.Range("O5").Value = 23.88141294715922
$ws.Range("B6").Value = 10.13122605778244
$ws.Range("C6").Value = 4.885676969033166
$ws.Range("D6").Value = 5.76429203539461
$ws.Range("E6").Value = 14.60985564899437
$ws.Range("G6").Value = 3.645470469524113
$ws.Range("K6").Value = 9.556339836512267
$ws.Range("N6").Value = 18.78687491563047
$ws.Range("O6").Value = 23.88317915273563
$ws.Range("B7").Value = 10.22213266799592
$ws.Range("C7").Value = 4.960219536614973
$ws.Range("D7").Value = 5.79782905918431
$ws.Range("E7").Value = 14.89536621600016
$ws.Range("G7").Value = 3.644760197126092
$ws.Range("K7").Value = 9.619051653475543
$ws.Range("N7").Value = 18.76656406277166
$ws.Range("O7").Value = 23.87143894386323
$ws.Range("B8").Value = 10.61929108111452
$ws.Range("C8").Value = 5.274003085567772
$ws.Range("D8").Value = 5.947390409104071
$ws.Range("E8").Value = 16.10337712160204
$ws.Range("G8").Value = 3.641785966673924
$ws.Range("K8").Value = 9.895851189292181
$ws.Range("N8").Value = 18.68095265417253
$ws.Range("O8").Value = 23.8325096420801
$ws.Range("B9").Value = 11.37556813204127
$ws.Range("C9").Value = 5.832244067627734
$ws.Range("D9").Value = 6.243385733432624
$ws.Range("E9").Value = 18.33009524098321
$ws.Range("G9").Value = 3.636524475497807
$ws.Range("K9").Value = 10.43364020119027
$ws.Range("N9").Value = 18.52740767558506
$ws.Range("O9").Value = 23.80270042369729
$ws.Range("B10").Value = 11.91091143197322
$ws.Range("C10").Value = 6.205820707942469
$ws.Range("D10").Value = 6.459968042753148
$ws.Range("E10").Value = 19.96063086427393
$ws.Range("G10").Value = 3.633002887455881
$ws.Range("K10").Value = 10.82128040108345
$ws.Range("N10").Value = 18.42324679391308
$ws.Range("O10").Value = 23.80934445215443
$ws.Range("B11").Value = 12.14884035552683
$ws.Range("C11").Value = 6.367560235273601
$ws.Range("D11").Value = 6.557798001531098
$ws.Range("E11").Value = 20.66055777284254
$ws.Range("G11").Value = 3.631474665516635
$ws.Range("K11").Value = 10.99514483028896
$ws.Range("N11").Value = 18.37771946446499
$ws.Range("O11").Value = 23.81858943217551
$ws.Range("B12").Value = 12.23804395462225
$ws.Range("C12").Value = 6.427609912225221
$ws.Range("D12").Value = 6.594704193982347
$ws.Range("E12").Value = 20.91962573156136
$ws.Range("G12").Value = 3.630906509306465
$ws.Range("K12").Value = 11.06056086429918
$ws.Range("N12").Value = 18.36074487821001
$ws.Range("O12").Value = 23.82298576685426
$ws.Range("B13").Value = 12.21887349742457
$ws.Range("C13").Value = 6.414730653777142
$ws.Range("D13").Value = 6.586762620778893
$ws.Range("E13").Value = 20.86409620815664
$ws.Range("G13").Value = 3.631028403714249
$ws.Range("K13").Value = 11.0464921635189
$ws.Range("N13").Value = 18.3643888652839
$ws.Range("O13").Value = 23.82199910901354
$ws.Range("B14").Value = 12.15619760635326
$ws.Range("C14").Value = 6.37252465082757
$ws.Range("D14").Value = 6.560837321297017
$ws.Range("E14").Value = 20.68199109108119
$ws.Range("G14").Value = 3.631427711897788
$ws.Range("K14").Value = 11.00053548372073
$ws.Range("N14").Value = 18.37631763870145
$ws.Range("O14").Value = 23.818933173759
$ws.Range("B15").Value = 12.11768776193974
$ws.Range("C15").Value = 6.346515818424216
$ws.Range("D15").Value = 6.544937958250408
$ws.Range("E15").Value = 20.56966876958821
$ws.Range("G15").Value = 3.631673671635982
$ws.Range("K15").Value = 10.97232874418182
$ws.Range("N15").Value = 18.38365891184637
$ws.Range("O15").Value = 23.81717182230673
$ws.Range("B16").Value = 11.89524165734534
$ws.Range("C16").Value = 6.19508397101002
$ws.Range("D16").Value = 6.453557179716373
$ws.Range("E16").Value = 19.91405048638637
$ws.Range("G16").Value = 3.633104239647492
$ws.Range("K16").Value = 10.8098623402033
$ws.Range("N16").Value = 18.42625935084289
$ws.Range("O16").Value = 23.80886555875679
$ws.Range("B17").Value = 11.75727710041045
$ws.Range("C17").Value = 6.100070617755867
$ws.Range("D17").Value = 6.397291007806162
$ws.Range("E17").Value = 19.50116465598972
$ws.Range("G17").Value = 3.634000697826578
$ws.Range("K17").Value = 10.70951164518054
$ws.Range("N17").Value = 18.45286779935276
$ws.Range("O17").Value = 23.80536439400871
$ws.Range("B18").Value = 11.67740108617956
$ws.Range("C18").Value = 6.04465024374057
$ws.Range("D18").Value = 6.364865094122393
$ws.Range("E18").Value = 19.25974613233914
$ws.Range("G18").Value = 3.634523263593471
$ws.Range("K18").Value = 10.65156332926889
$ws.Range("N18").Value = 18.46834704321716
$ws.Range("O18").Value = 23.80393632227413
$ws.Range("B19").Value = 11.65026957739346
$ws.Range("C19").Value = 6.025754010884552
$ws.Range("D19").Value = 6.353876654637471
$ws.Range("E19").Value = 19.17732867701483
$ws.Range("G19").Value = 3.63470139030109
$ws.Range("K19").Value = 10.63190575796282
$ws.Range("N19").Value = 18.47361810448174
$ws.Range("O19").Value = 23.80355336495277
$ws.Range("B20").Value = 11.77201841015338
$ws.Range("C20").Value = 6.110264910416403
$ws.Range("D20").Value = 6.403287465363272
$ws.Range("E20").Value = 19.54552433338482
$ws.Range("G20").Value = 3.633904549792203
$ws.Range("K20").Value = 10.72021837216592
$ws.Range("N20").Value = 18.45001720588833
$ws.Range("O20").Value = 23.80567647340084
$ws.Range("B21").Value = 12.17463199256285
$ws.Range("C21").Value = 6.38495420581944
$ws.Range("D21").Value = 6.568456314924084
$ws.Range("E21").Value = 20.73564178670197
$ws.Range("G21").Value = 3.631310139598249
$ws.Range("K21").Value = 11.0140460549268
$ws.Range("N21").Value = 18.37280666979584
$ws.Range("O21").Value = 23.81980940911562
$ws.Range("B22").Value = 12.43251151524993
$ws.Range("C22").Value = 6.557493112766986
$ws.Range("D22").Value = 6.675571326723917
$ws.Range("E22").Value = 21.47861646628242
$ws.Range("G22").Value = 3.629675998724234
$ws.Range("K22").Value = 11.20358811045951
$ws.Range("N22").Value = 18.32389281115453
$ws.Range("O22").Value = 23.8342653125746
$ws.Range("B23").Value = 12.29538421955277
$ws.Range("C23").Value = 6.466050216855485
$ws.Range("D23").Value = 6.61849084468925
$ws.Range("E23").Value = 21.08525342603632
$ws.Range("G23").Value = 3.630542566695736
$ws.Range("K23").Value = 11.10267479290721
$ws.Range("N23").Value = 18.34985785844847
$ws.Range("O23").Value = 23.82607234104321
$ws.Range("B24").Value = 11.76535559624247
$ws.Range("C24").Value = 6.10565854853566
$ws.Range("D24").Value = 6.400576706006708
$ws.Range("E24").Value = 19.52548192562957
$ws.Range("G24").Value = 3.633947995947665
$ws.Range("K24").Value = 10.71537865113822
$ws.Range("N24").Value = 18.45130539307069
$ws.Range("O24").Value = 23.80553356064346
$ws.Range("B25").Value = 11.1740996332629
$ws.Range("C25").Value = 5.687539166974092
$ws.Range("D25").Value = 6.163288642761314
$ws.Range("E25").Value = 17.70689624268311
$ws.Range("G25").Value = 3.637887142036106
$ws.Range("K25").Value = 10.28916076479682
$ws.Range("N25").Value = 18.56742021949639
$ws.Range("O25").Value = 23.82199910901354
